# Daily attendance processing - 2025-12-14 06:36:38
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("L6").Value = 59

$ws.Range("L7").Value = 23

$ws.Range("G8").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("H8").Value = "105/216"
$ws.Range("L8").Value = 94

$ws.Range("L9").Value = "33.5%"

$ws.Range("G10").Value = "arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("L10").Value = "35.5%"

$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 11
$ws.Range("S15").Value = "39.8%"

$ws.Range("O16").Value = 9
$ws.Range("P16").Value = 1
$ws.Range("R16").Value = "40.9%"
$ws.Range("S16").Value = "40.6%"

$ws.Range("P17").Value = 4
$ws.Range("Q17").Value = 11

$ws.Range("G19").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G20").Value = "youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("A3:I3").Copy()
$ws.Range("A21:I21").PasteSpecial(-4122)
$ws.Range("I21").Value = "Not Recorded"
$ws.Range("S21").Value = "30.2%"

$ws.Range("S22").Value = "19.4%"

$ws.Range("G24").Value = "nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("G25").Value = "manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, backup@backdoor.com, gehanadel@med.asu.edu.eg"

$ws.Range("A2:I2").Copy()
$ws.Range("A30:I30").PasteSpecial(-4122)
$ws.Range("G30").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("H30").Value = "106/217"
$ws.Range("I30").Value = "Recorded"

$ws.Range("G32").Value = "arwaelsayed03@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"

$ws.Range("G41").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G42").Value = "youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G46").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"

$ws.Range("G54").Value = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("A3:I3").Copy()
$ws.Range("A55:I55").PasteSpecial(-4122)
$ws.Range("I55").Value = "Not Recorded"

$ws.Range("G62").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G63").Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

$ws.Range("G64").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G68").Value = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"

$ws.Range("G76").Value = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("G81").Value = "enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, user@user.com"

$ws.Range("G84").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G85").Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

$ws.Range("G86").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G98").Value = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("G106").Value = "youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G108").Value = "youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G120").Value = "maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("G128").Value = "youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G130").Value = "youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G134").Value = "hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("G140").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("H140").Value = "24/224"

$ws.Range("G142").Value = "basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("G150").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G156").Value = "manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("G162").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("H162").Value = "40/226"

$ws.Range("G164").Value = "basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$ws.Range("G172").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$excel.CutCopyMode = 0
